# flag_map.xlsx: add a new flag-map entry for missing administration term units.
# A new row is inserted at row 18 (pushing the existing rows 18-54 down to 19-55),
# populated with the new Record Identifier / Definition / Flag Type, and the
# sheet's AutoFilter / _FilterDatabase range is extended to cover the extra row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18, shifting existing data down.
$ws.Rows.Item(18).Insert()

# Populate the new row with the new flag-map entry.
$ws.Range("A18").Value = "missing_administration_term_units"
$ws.Range("B18").Value = "Missing administration term units for conversion"
$ws.Range("C18").Value = "Hard Stop (Missing Required)"

# Re-apply AutoFilter so its range grows from A1:C54 to A1:C55.
$ws.AutoFilterMode = $false
$ws.Range("A1:C55").AutoFilter()

# Keep the _FilterDatabase defined name in sync with the new range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$C`$55"
    }
}

# Match the recorded selection/cursor position after the edit.
$ws.Range("A18").Select()
